$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "2020-12-11"
$ws.Range("J2").Value = 2500
$ws.Range("K2").Value = 1200
$ws.Range("M2").Value = 1256
$ws.Range("P2").Value = 1256
$ws.Range("D3").Value = "2020-12-11"
$ws.Range("J3").Value = 1500
$ws.Range("D4").Value = "2020-12-16"
$ws.Range("J4").Value = 1000
$ws.Range("D5").Value = "2020-12-16"
$ws.Range("J5").Value = 900
$ws.Range("K5").Value = 900
$ws.Range("L5").Value = 900
$ws.Range("M5").Value = 900
$ws.Range("P5").Value = 900
$ws.Range("D6").Value = "2020-12-09"
$ws.Range("J6").Value = 2800
$ws.Range("L6").Value = 1250
$ws.Range("M6").Value = 1221
$ws.Range("O6").Value = "Provincia de Quillota"
$ws.Range("P6").Value = 1221
$ws.Range("D7").Value = "2020-12-09"
$ws.Range("I7").Value = "Segunda"
$ws.Range("J7").Value = 1300
$ws.Range("K7").Value = 1000
$ws.Range("L7").Value = 1000
$ws.Range("M7").Value = 1000
$ws.Range("P7").Value = 1000
$ws.Range("D8").Value = "2020-11-30"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 720
$ws.Range("K8").Value = 1200
$ws.Range("L8").Value = 1200
$ws.Range("M8").Value = 1200
$ws.Range("P8").Value = 1200
$ws.Range("D9").Value = "2020-11-30"
$ws.Range("I9").Value = "Segunda"
$ws.Range("J9").Value = 750
$ws.Range("K9").Value = 1000
$ws.Range("L9").Value = 1000
$ws.Range("M9").Value = 1000
$ws.Range("P9").Value = 1000
$ws.Range("D10").Value = "2020-11-24"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 1100
$ws.Range("K10").Value = 1300
$ws.Range("L10").Value = 1300
$ws.Range("M10").Value = 1300
$ws.Range("P10").Value = 1300
$ws.Range("D11").Value = "2020-11-24"
$ws.Range("I11").Value = "Segunda"
$ws.Range("J11").Value = 800
$ws.Range("K11").Value = 1000
$ws.Range("L11").Value = 1000
$ws.Range("M11").Value = 1000
$ws.Range("P11").Value = 1000
$ws.Range("D12").Value = "2020-12-10"
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 1500
$ws.Range("K12").Value = 1300
$ws.Range("L12").Value = 1300
$ws.Range("M12").Value = 1300
$ws.Range("P12").Value = 1300
$ws.Range("D13").Value = "2020-12-10"
$ws.Range("I13").Value = "Segunda"
$ws.Range("J13").Value = 1450
$ws.Range("K13").Value = 1000
$ws.Range("L13").Value = 1000
$ws.Range("M13").Value = 1000
$ws.Range("P13").Value = 1000
$ws.Range("D14").Value = "2020-12-07"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 600
$ws.Range("K14").Value = 1300
$ws.Range("L14").Value = 1300
$ws.Range("M14").Value = 1300
$ws.Range("P14").Value = 1300
$ws.Range("D15").Value = "2020-12-07"
$ws.Range("I15").Value = "Segunda"
$ws.Range("J15").Value = 550
$ws.Range("K15").Value = 1000
$ws.Range("L15").Value = 1000
$ws.Range("M15").Value = 1000
$ws.Range("P15").Value = 1000
$ws.Range("D16").Value = "2020-12-03"
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 1200
$ws.Range("K16").Value = 1300
$ws.Range("L16").Value = 1300
$ws.Range("M16").Value = 1300
$ws.Range("P16").Value = 1300
$ws.Range("D17").Value = "2020-12-03"
$ws.Range("I17").Value = "Segunda"
$ws.Range("J17").Value = 850
$ws.Range("K17").Value = 1000
$ws.Range("L17").Value = 1000
$ws.Range("M17").Value = 1000
$ws.Range("P17").Value = 1000
$ws.Range("D18").Value = "2020-12-14"
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 980
$ws.Range("K18").Value = 1200
$ws.Range("L18").Value = 1200
$ws.Range("M18").Value = 1200
$ws.Range("O18").Value = "Región Metropolitana"
$ws.Range("P18").Value = 1200
$ws.Range("D19").Value = "2020-11-26"
$ws.Range("J19").Value = 1600
$ws.Range("D20").Value = "2020-11-26"
$ws.Range("J20").Value = 1850
$ws.Range("D21").Value = "2020-11-25"
$ws.Range("J21").Value = 750
$ws.Range("K21").Value = 1300
$ws.Range("M21").Value = 1300
$ws.Range("P21").Value = 1300
$ws.Range("D22").Value = "2020-11-25"
$ws.Range("J22").Value = 850
$ws.Range("D23").Value = "2020-12-02"
$ws.Range("J23").Value = 1430
$ws.Range("K23").Value = 1200
$ws.Range("M23").Value = 1248
$ws.Range("P23").Value = 1248
$ws.Range("D24").Value = "2020-12-02"
$ws.Range("J24").Value = 350
$ws.Range("D25").Value = "2020-12-04"
$ws.Range("J25").Value = 950
$ws.Range("D26").Value = "2020-12-04"
$ws.Range("J26").Value = 800
$ws.Range("D27").Value = "2020-11-27"
$ws.Range("J27").Value = 1200
$ws.Range("D28").Value = "2020-11-27"
$ws.Range("J28").Value = 800
